$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date (2022-02-23 -> 2022-02-24)
$ws.Name = "Through 2022-02-24"

# Update the February row label to match the new "through" date
$ws.Range("A3").Value = "February (through 02-24)"

# Update February row (row 3) year totals
$ws.Range("C3").Value = 30
$ws.Range("D3").Value = 51
$ws.Range("G3").Value = 61
$ws.Range("H3").Value = 108
$ws.Range("I3").Value = 121

# Update overall Total row (row 4) year totals
$ws.Range("C4").Value = 81
$ws.Range("D4").Value = 126
$ws.Range("G4").Value = 135
$ws.Range("H4").Value = 325
$ws.Range("I4").Value = 280
